# Append new price/NAV history rows (2024-08-28 .. 2024-09-25) to the
# Pharma_stocks tracking sheet, mirroring the existing row layout
# (columns A, C-J; column B is unused in data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 618; A = "2024-08-28"; C = 2200.75; D = 1539.5; E = 1707.449951171875; F = 1969.050048828125; G = 1138.300048828125; H = 8555.050048828125; I = 0.0; J = 189.512658754004 },
    @{ Row = 619; A = "2024-08-29"; C = 2193.75; D = 1499.150024414062; E = 1691.300048828125; F = 1961.150024414062; G = 1132.050048828125; H = 8477.400146484375; I = -0.009076498898377167; J = 187.7925473155948 },
    @{ Row = 620; A = "2024-08-30"; C = 2240.199951171875; D = 1537.550048828125; E = 1731.75; F = 1953.800048828125; G = 1127.900024414062; H = 8591.200073242188; I = 0.01342391827581785; J = 190.313459123567 },
    @{ Row = 621; A = "2024-09-02"; C = 2232.75; D = 1537.550048828125; E = 1687.900024414062; F = 1970.599975585938; G = 1111.550048828125; H = 8540.35009765625; I = -0.005918844300264037; J = 189.1870233907699 },
    @{ Row = 622; A = "2024-09-03"; C = 2240.25; D = 1530.599975585938; E = 1687.5; F = 1924.650024414062; G = 1114.0; H = 8497.0; I = -0.005075915759957742; J = 188.2267259971612 },
    @{ Row = 623; A = "2024-09-04"; C = 2277.25; D = 1556.550048828125; E = 1686.550048828125; F = 1924.650024414062; G = 1127.900024414062; H = 8572.900146484375; I = 0.008932581674046723; J = 189.9080766003692 },
    @{ Row = 624; A = "2024-09-05"; C = 2290.199951171875; D = 1555.75; E = 1709.449951171875; F = 1933.599975585938; G = 1115.150024414062; H = 8604.14990234375; I = 0.00364517903223101; J = 190.6003255392442 },
    @{ Row = 625; A = "2024-09-06"; C = 2256.5; D = 1559.900024414062; E = 1702.699951171875; F = 1928.400024414062; G = 1100.0; H = 8547.5; I = -0.006584020848860234; J = 189.3454090220943 },
    @{ Row = 626; A = "2024-09-09"; C = 2216.800048828125; D = 1546.25; E = 1704.199951171875; F = 1937.099975585938; G = 1104.150024414062; H = 8508.5; I = -0.004562737642585551; J = 188.4814755968984 },
    @{ Row = 627; A = "2024-09-10"; C = 2222.550048828125; D = 1545.550048828125; E = 1727.849975585938; F = 1912.150024414062; G = 1113.199951171875; H = 8521.300048828125; I = 0.001504383713712758; J = 188.765024059123 },
    @{ Row = 628; A = "2024-09-11"; C = 2209.39990234375; D = 1591.949951171875; E = 1725.650024414062; F = 1867.75; G = 1112.599975585938; H = 8507.349853515625; I = -0.001637097066476197; J = 188.4559973919825 },
    @{ Row = 629; A = "2024-09-12"; C = 2247.5; D = 1592.849975585938; E = 1747.949951171875; F = 1883.349975585938; G = 1120.099975585938; H = 8591.749877929688; I = 0.009920836202496664; J = 190.3256384734865 },
    @{ Row = 630; A = "2024-09-13"; C = 2256.449951171875; D = 1582.5; E = 1753.699951171875; F = 1923.300048828125; G = 1118.550048828125; H = 8634.5; I = 0.004975717714982386; J = 191.2726451244543 },
    @{ Row = 631; A = "2024-09-16"; C = 2251.85009765625; D = 1577.75; E = 1741.449951171875; F = 1900.949951171875; G = 1115.849975585938; H = 8587.849975585938; I = -0.005402747630327465; J = 190.2392472942617 },
    @{ Row = 632; A = "2024-09-17"; C = 2270.39990234375; D = 1561.699951171875; E = 1713.0; F = 1875.599975585938; G = 1110.949951171875; H = 8531.649780273438; I = -0.006544151967287428; J = 188.9942927498257 },
    @{ Row = 633; A = "2024-09-18"; C = 2224.949951171875; D = 1543.050048828125; E = 1646.050048828125; F = 1857.0; G = 1079.949951171875; H = 8351.0; I = -0.02117407358786916; J = 184.9925136874536 },
    @{ Row = 634; A = "2024-09-19"; C = 2171.89990234375; D = 1515.050048828125; E = 1649.800048828125; F = 1886.5; G = 1054.449951171875; H = 8277.699951171875; I = -0.00877739777608969; J = 183.3687608092201 },
    @{ Row = 635; A = "2024-09-20"; C = 2151.699951171875; D = 1481.099975585938; E = 1636.75; F = 1897.25; G = 1054.599975585938; H = 8221.39990234375; I = -0.0068014121265841; J = 182.1215942958156 },
    @{ Row = 636; A = "2024-09-23"; C = 2182.25; D = 1440.400024414062; E = 1712.449951171875; F = 1952.0; G = 1055.25; H = 8342.349975585938; I = 0.01471161537923817; J = 184.8008971433493 },
    @{ Row = 637; A = "2024-09-24"; C = 2215.75; D = 1414.25; E = 1697.5; F = 1944.349975585938; G = 1051.550048828125; H = 8323.400024414062; I = -0.002271536344954651; J = 184.3811151889079 },
    @{ Row = 638; A = "2024-09-25"; C = 2221.10009765625; D = 1416.400024414062; E = 1689.199951171875; F = 1909.550048828125; G = 1063.449951171875; H = 8299.700073242188; I = -0.002847388219040138; J = 183.8561105737055 }
)

$firstRow = $newRows[0].Row
$lastRow = $newRows[$newRows.Count - 1].Row

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Column A holds the date as literal text (e.g. "2024-08-28"), matching
    # the existing rows above it. Assigning a date-shaped string straight to
    # .Value would get auto-converted to a date serial by Excel's input
    # parser, so we instead enter it as a text formula and then flatten the
    # formula down to its cached string value via Copy / PasteSpecial values
    # (this sidesteps the "smart" literal parsing entirely and leaves the
    # cell as a plain text value like the original, with no formula and no
    # number-format override left behind).
    $ws.Cells.Item($rowNum, 1).Formula = '="' + $r.A + '"'

    $ws.Cells.Item($rowNum, 3).Value  = $r.C
    $ws.Cells.Item($rowNum, 4).Value  = $r.D
    $ws.Cells.Item($rowNum, 5).Value  = $r.E
    $ws.Cells.Item($rowNum, 6).Value  = $r.F
    $ws.Cells.Item($rowNum, 7).Value  = $r.G
    $ws.Cells.Item($rowNum, 8).Value  = $r.H
    $ws.Cells.Item($rowNum, 9).Value  = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
}

$dateRange = $ws.Range("A" + $firstRow + ":A" + $lastRow)
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0
